$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Step 1: locate " lors de API code." (plus the trailing lone-space run
# right after it) at the end of the last paragraph, and replace that whole
# span with the reworded / re-split runs. ---
$find = $d.Content
$found = $find.Find.Execute(" lors de API code.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence ' lors de API code.'"
}

$targetStart = $find.Start

# Grow the end to cover the trailing run (a lone space) that sits right
# before the paragraph mark, so it gets folded into the replacement too.
$paraRange = $find.Duplicate
$paraRange.Expand(4)  # wdParagraph
$paraStart = $paraRange.Start
$targetEnd = $paraRange.End - 1

$replaceRange = $d.Range($targetStart, $targetEnd)

$newRunsXml = '<w:p ' + $wNs + '>' `
    + '<w:r w:rsidR="00904438"><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:proofErr w:type="gramStart"/>' `
    + '<w:r><w:t>lors de API</w:t></w:r>' `
    + '<w:proofErr w:type="gramEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> code.</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> En fait non, c’est un problème de cache qui recharge souvent.</w:t></w:r>' `
    + '</w:p>'

$replaceRange.InsertXML($newRunsXml)

# --- Step 2: append three brand-new paragraphs right after that paragraph. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

$newParasXml = '<w:p ' + $wNs + '>' `
        + '<w:r><w:t>Les caches des navigateurs peuvent vraiment être agaçants quand on fait du développement. La prochaine fois, pense à utiliser Ctrl + Shift + R (ou Cmd + Shift + R sur Mac) pour forcer un rechargement sans cache.</w:t></w:r>' `
    + '</w:p>' `
    + '<w:p ' + $wNs + '>' `
        + '<w:r><w:t>Si le problème revient souvent, tu peux aussi :</w:t></w:r>' `
    + '</w:p>' `
    + '<w:p ' + $wNs + '>' `
        + '<w:r><w:t xml:space="preserve"> Désactiver le cache dans les </w:t></w:r>' `
        + '<w:proofErr w:type="spellStart"/>' `
        + '<w:r><w:t>DevTools</w:t></w:r>' `
        + '<w:proofErr w:type="spellEnd"/>' `
        + '<w:r><w:t xml:space="preserve"> (F12 &gt; onglet "Network" &gt; cocher "</w:t></w:r>' `
        + '<w:proofErr w:type="spellStart"/>' `
        + '<w:r><w:t>Disable</w:t></w:r>' `
        + '<w:proofErr w:type="spellEnd"/>' `
        + '<w:r><w:t xml:space="preserve"> cache")</w:t></w:r>' `
    + '</w:p>'

$insertPoint.InsertXML($newParasXml)

Write-Host "Edit applied."
